$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B142").Value = 48654
$ws.Range("E142").Value = 38.26
$ws.Range("F142").Value = -1
$ws.Range("G142").Value = -32.02

$ws.Range("B143").Value = 63902
$ws.Range("E143").Value = 34.04
$ws.Range("F143").Value = 2
$ws.Range("G143").Value = 64.04000000000001

$ws.Range("B256").Value = 64979
$ws.Range("E256").Value = 314.41
$ws.Range("F256").Value = 82
$ws.Range("G256").Value = 24251.5

$ws.Range("B257").Value = 48719
$ws.Range("E257").Value = 353.35
$ws.Range("F257").Value = -81
$ws.Range("G257").Value = -23955.75

$ws.Range("B271").Value = 48706
$ws.Range("E271").Value = 39.8
$ws.Range("F271").Value = -144
$ws.Range("G271").Value = -4795.2

$ws.Range("B272").Value = 64973
$ws.Range("E272").Value = 35.4
$ws.Range("F272").Value = 150
$ws.Range("G272").Value = 4995

$ws.Range("B305").Value = 62997
$ws.Range("F305").Value = 72
$ws.Range("G305").Value = 22020.48

$ws.Range("B306").Value = 57854
$ws.Range("F306").Value = 2
$ws.Range("G306").Value = 611.6799999999999

$ws.Range("B309").Value = 61610
$ws.Range("D309").Value = 102.71
$ws.Range("E309").Value = 122.71
$ws.Range("F309").Value = -58
$ws.Range("G309").Value = -5957.18

$ws.Range("B310").Value = 57077
$ws.Range("D310").Value = 93.08
$ws.Range("E310").Value = 111.2
$ws.Range("F310").Value = 1
$ws.Range("G310").Value = 93.08

$ws.Range("B338").Value = 63520
$ws.Range("E338").Value = 153.4
$ws.Range("F338").Value = 97
$ws.Range("G338").Value = 13995.16

$ws.Range("B339").Value = 55373
$ws.Range("E339").Value = 163.62
$ws.Range("F339").Value = -94
$ws.Range("G339").Value = -13562.32

$ws.Range("B342").Value = 63571
$ws.Range("E342").Value = 152.53
$ws.Range("F342").Value = 29
$ws.Range("G342").Value = 4160.92

$ws.Range("B343").Value = 63531
$ws.Range("F343").Value = 80
$ws.Range("G343").Value = 11478.4

$ws.Range("B344").Value = 57802
$ws.Range("E344").Value = 162.71
$ws.Range("F344").Value = -79
$ws.Range("G344").Value = -11334.92

$ws.Range("B364").Value = 63652
$ws.Range("E364").Value = 55.42
$ws.Range("F364").Value = 250
$ws.Range("G364").Value = 13032.5

$ws.Range("B365").Value = 57885
$ws.Range("E365").Value = 62.28
$ws.Range("F365").Value = 4
$ws.Range("G365").Value = 208.52

$ws.Range("B371").Value = 61608
$ws.Range("E371").Value = 154.12
$ws.Range("F371").Value = -56
$ws.Range("G371").Value = -7224.56

$ws.Range("B372").Value = 63564
$ws.Range("E372").Value = 137.16
$ws.Range("F372").Value = 57
$ws.Range("G372").Value = 7353.57

$ws.Range("B381").Value = 62865
$ws.Range("F381").Value = 151
$ws.Range("G381").Value = 12051.31

$ws.Range("B382").Value = 57817
$ws.Range("F382").Value = 3
$ws.Range("G382").Value = 239.43

$ws.Range("B392").Value = 57835
$ws.Range("F392").Value = 1
$ws.Range("G392").Value = 59.13

$ws.Range("B393").Value = 62933
$ws.Range("F393").Value = 146
$ws.Range("G393").Value = 8632.98

$ws.Range("B411").Value = 63007
$ws.Range("F411").Value = 984
$ws.Range("G411").Value = 168588.72

$ws.Range("B412").Value = 57856
$ws.Range("F412").Value = 2
$ws.Range("G412").Value = 342.66

$ws.Range("B413").Value = 63008
$ws.Range("F413").Value = 504
$ws.Range("G413").Value = 76189.67999999999

$ws.Range("B414").Value = 57857
$ws.Range("F414").Value = 3
$ws.Range("G414").Value = 453.51

$ws.Range("B423").Value = 53082
$ws.Range("C423").Value = 'HUL-VIM BAR MULTIPACK FW 4X200G'
$ws.Range("F423").Value = 1
$ws.Range("G423").Value = 59.47

$ws.Range("B424").Value = 63102
$ws.Range("C424").Value = 'HUL-Vim Bar Multipack Fw 4X200G'
$ws.Range("F424").Value = 36
$ws.Range("G424").Value = 2140.92

$ws.Range("B449").Value = 31930
$ws.Range("E449").Value = 26.8
$ws.Range("F449").Value = -62
$ws.Range("G449").Value = -1390.04

$ws.Range("B450").Value = 63681
$ws.Range("E450").Value = 23.84
$ws.Range("F450").Value = 65
$ws.Range("G450").Value = 1457.3

$ws.Range("B528").Value = 58047
$ws.Range("D528").Value = 105.54
$ws.Range("E528").Value = 126.1
$ws.Range("F528").Value = 54
$ws.Range("G528").Value = 5699.16

$ws.Range("B529").Value = 47097
$ws.Range("D529").Value = 112.28
$ws.Range("E529").Value = 134.16
$ws.Range("F529").Value = 15
$ws.Range("G529").Value = 1684.2

$ws.Range("B571").Value = 53757
$ws.Range("E571").Value = 16.08
$ws.Range("F571").Value = -159
$ws.Range("G571").Value = -2138.55

$ws.Range("B572").Value = 65069
$ws.Range("E572").Value = 14.3
$ws.Range("F572").Value = 172
$ws.Range("G572").Value = 2313.4

$ws.Range("B575").Value = 65066
$ws.Range("E575").Value = 13.61
$ws.Range("F575").Value = 313
$ws.Range("G575").Value = 4009.53

$ws.Range("B576").Value = 53263
$ws.Range("E576").Value = 15.29
$ws.Range("F576").Value = -309
$ws.Range("G576").Value = -3958.29

$ws.Range("B578").Value = 45695
$ws.Range("E578").Value = 23.58
$ws.Range("F578").Value = -36
$ws.Range("G578").Value = -710.28

$ws.Range("B579").Value = 64915
$ws.Range("E579").Value = 20.98
$ws.Range("F579").Value = 40
$ws.Range("G579").Value = 789.2

$ws.Range("B585").Value = 64927
$ws.Range("E585").Value = 17.26
$ws.Range("F585").Value = 295
$ws.Range("G585").Value = 4784.9

$ws.Range("B586").Value = 45718
$ws.Range("E586").Value = 19.38
$ws.Range("F586").Value = -294
$ws.Range("G586").Value = -4768.68

$ws.Range("B591").Value = 45709
$ws.Range("E591").Value = 15.69
$ws.Range("F591").Value = -300
$ws.Range("G591").Value = -3945

$ws.Range("B592").Value = 64925
$ws.Range("E592").Value = 13.97
$ws.Range("F592").Value = 302
$ws.Range("G592").Value = 3971.3

$ws.Range("B679").Value = 53319
$ws.Range("E679").Value = 310.64
$ws.Range("F679").Value = -6
$ws.Range("G679").Value = -1643.52

$ws.Range("B680").Value = 64810
$ws.Range("E680").Value = 291.22
$ws.Range("F680").Value = 7
$ws.Range("G680").Value = 1917.44

$ws.Range("B701").Value = 64833
$ws.Range("E701").Value = 34.9
$ws.Range("F701").Value = 99
$ws.Range("G701").Value = 3250.17

$ws.Range("B702").Value = 60025
$ws.Range("E702").Value = 37.22
$ws.Range("F702").Value = -98
$ws.Range("G702").Value = -3217.34

$ws.Range("B707").Value = 60031
$ws.Range("E707").Value = 111.69
$ws.Range("F707").Value = -5
$ws.Range("G707").Value = -492.5

$ws.Range("B708").Value = 64836
$ws.Range("E708").Value = 104.71
$ws.Range("F708").Value = 7
$ws.Range("G708").Value = 689.5

$ws.Range("B712").Value = 64830
$ws.Range("E712").Value = 34.9
$ws.Range("F712").Value = 117
$ws.Range("G712").Value = 3841.11

$ws.Range("B713").Value = 60022
$ws.Range("E713").Value = 37.22
$ws.Range("F713").Value = -113
$ws.Range("G713").Value = -3709.79

$ws.Range("B864").Value = 65079
$ws.Range("E864").Value = 43.44
$ws.Range("F864").Value = 21
$ws.Range("G864").Value = 858.27

$ws.Range("B865").Value = 54751
$ws.Range("E865").Value = 46.34
$ws.Range("F865").Value = -19
$ws.Range("G865").Value = -776.53
